$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new row above row 5 (current top data row), pulling formatting
# from the row below so the header row's bold/fill style isn't inherited.
$ws.Rows.Item(5).Insert(-4121)
$ws.Rows.Item(5).ClearFormats()

# New row 5 data (most recent trade), matching style of the date column (s=2)
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A5").Value = 46066
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 2003.6
$ws.Range("F5").Value = 6053.16
$ws.Range("G5").Value = "CN#252611910666"
$ws.Range("H5").Value = 6.0717
$ws.Range("I5").Value = 36.2857
$ws.Range("J5").Formula = '=Index!$C$2'
